# Insert a new data row at row 12 (pushing the existing rows 12-80 down
# to 13-81) and populate it with the new "Murcott / Primera" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 44532
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100102
$ws.Range("H12").Value = "Cítricos"
$ws.Range("I12").Value = 100102004
$ws.Range("J12").Value = "Mandarina"
$ws.Range("K12").Value = "Murcott"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 270
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14500
$ws.Range("Q12").Value = "$/caja 20 kilos"
$ws.Range("R12").Value = "Región de Coquimbo"
$ws.Range("S12").Value = 725
$ws.Range("T12").Value = 20
